# Update automàtic: dades i banners [2026-02-20 07:26]
# Applies the latest meteocat extraction pass: refreshed DATA_EXTRACCIO
# timestamps and the handful of station readings (humitat, pressio,
# radiacio, vent, temperatures) that shifted between pulls.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-20 07:24:32'
$ws.Range('E3').Value = '2026-02-20 07:24:34'
$ws.Range('N3').Value = '-6.3 °C 6:31 TU'
$ws.Range('E4').Value = '2026-02-20 07:24:37'
$ws.Range('J4').Value = '1019.1 hPa'
$ws.Range('O4').Value = '8.2 °C'
$ws.Range('E5').Value = '2026-02-20 07:24:39'
$ws.Range('N5').Value = '-6.0 °C 6:59 TU'
$ws.Range('E6').Value = '2026-02-20 07:24:42'
$ws.Range('J6').Value = '1019.2 hPa'
$ws.Range('L6').Value = '15.5 km/h - 285º 6:47 TU'
$ws.Range('E7').Value = '2026-02-20 07:24:44'
$ws.Range('H7').NumberFormat = '@'
$ws.Range('H7').Value = '47%'
$ws.Range('H7').NumberFormat = 'general'
$ws.Range('J7').Value = '1018.7 hPa'
$ws.Range('K7').Value = '0.0 MJ/m2'
$ws.Range('M7').Value = '11.7 °C 6:40 TU'
$ws.Range('O7').Value = '11.1 °C'
$ws.Range('E8').Value = '2026-02-20 07:24:47'
$ws.Range('H8').NumberFormat = '@'
$ws.Range('H8').Value = '60%'
$ws.Range('H8').NumberFormat = 'general'
$ws.Range('J8').Value = '1019.6 hPa'
$ws.Range('E9').Value = '2026-02-20 07:24:49'
$ws.Range('H9').NumberFormat = '@'
$ws.Range('H9').Value = '46%'
$ws.Range('H9').NumberFormat = 'general'
$ws.Range('L9').Value = '57.6 km/h - 6º 6:51 TU'
$ws.Range('O9').Value = '12.3 °C'
$ws.Range('E10').Value = '2026-02-20 07:24:52'
$ws.Range('E11').Value = '2026-02-20 07:24:54'
$ws.Range('O11').Value = '7.7 °C'
$ws.Range('E12').Value = '2026-02-20 07:24:56'
$ws.Range('E13').Value = '2026-02-20 07:24:59'
$ws.Range('J13').Value = '1020.7 hPa'
$ws.Range('N13').Value = '1.7 °C 6:32 TU'
$ws.Range('O13').Value = '4.7 °C'
$ws.Range('E14').Value = '2026-02-20 07:25:01'
$ws.Range('N14').Value = '9.3 °C 6:46 TU'
$ws.Range('E15').Value = '2026-02-20 07:25:04'
$ws.Range('H15').NumberFormat = '@'
$ws.Range('H15').Value = '45%'
$ws.Range('H15').NumberFormat = 'general'
$ws.Range('N15').Value = '11.1 °C 6:38 TU'
$ws.Range('O15').Value = '12.5 °C'
$ws.Range('E16').Value = '2026-02-20 07:25:06'
$ws.Range('E17').Value = '2026-02-20 07:25:08'
$ws.Range('K17').Value = '0.0 MJ/m2'
$ws.Range('E18').Value = '2026-02-20 07:25:11'
$ws.Range('J18').Value = '1019.6 hPa'
$ws.Range('N18').Value = '-0.4 °C 6:44 TU'
$ws.Range('O18').Value = '1.5 °C'
$ws.Range('E19').Value = '2026-02-20 07:25:13'
$ws.Range('E20').Value = '2026-02-20 07:25:16'
$ws.Range('H20').NumberFormat = '@'
$ws.Range('H20').Value = '73%'
$ws.Range('H20').NumberFormat = 'general'
$ws.Range('M20').Value = '-4.3 °C 6:59 TU'
$ws.Range('E21').Value = '2026-02-20 07:25:18'
$ws.Range('J21').Value = '1020.7 hPa'
$ws.Range('O21').Value = '5.9 °C'
$ws.Range('E22').Value = '2026-02-20 07:25:21'
$ws.Range('H22').NumberFormat = '@'
$ws.Range('H22').Value = '59%'
$ws.Range('H22').NumberFormat = 'general'
$ws.Range('M22').Value = '-3.7 °C 6:59 TU'
$ws.Range('O22').Value = '-6.3 °C'
$ws.Range('E23').Value = '2026-02-20 07:25:23'
$ws.Range('H23').NumberFormat = '@'
$ws.Range('H23').Value = '85%'
$ws.Range('H23').NumberFormat = 'general'
$ws.Range('I23').Value = '3.8 mm'
$ws.Range('N23').Value = '-7.1 °C 6:38 TU'
$ws.Range('O23').Value = '-6.6 °C'
$ws.Range('E24').Value = '2026-02-20 07:25:26'
$ws.Range('J24').Value = '1023.3 hPa'
$ws.Range('E25').Value = '2026-02-20 07:25:28'
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H25').Value = '68%'
$ws.Range('H25').NumberFormat = 'general'
$ws.Range('M25').Value = '-3.1 °C 6:53 TU'
$ws.Range('O25').Value = '-4.7 °C'
$ws.Range('E26').Value = '2026-02-20 07:25:30'
$ws.Range('J26').Value = '1018.9 hPa'
$ws.Range('E27').Value = '2026-02-20 07:25:33'
$ws.Range('H27').NumberFormat = '@'
$ws.Range('H27').Value = '50%'
$ws.Range('H27').NumberFormat = 'general'
$ws.Range('K27').Value = '0.0 MJ/m2'
$ws.Range('E28').Value = '2026-02-20 07:25:35'
$ws.Range('E29').Value = '2026-02-20 07:25:38'
$ws.Range('E30').Value = '2026-02-20 07:25:40'
$ws.Range('E31').Value = '2026-02-20 07:25:42'
$ws.Range('E32').Value = '2026-02-20 07:25:45'
$ws.Range('E33').Value = '2026-02-20 07:25:47'
$ws.Range('E34').Value = '2026-02-20 07:25:49'
$ws.Range('E35').Value = '2026-02-20 07:25:52'
$ws.Range('E36').Value = '2026-02-20 07:25:54'
$ws.Range('E37').Value = '2026-02-20 07:25:56'
$ws.Range('E38').Value = '2026-02-20 07:25:59'
$ws.Range('E39').Value = '2026-02-20 07:26:01'
$ws.Range('E40').Value = '2026-02-20 07:26:04'
$ws.Range('E41').Value = '2026-02-20 07:26:06'
$ws.Range('E42').Value = '2026-02-20 07:26:08'
$ws.Range('E43').Value = '2026-02-20 07:26:11'
$ws.Range('E44').Value = '2026-02-20 07:26:13'
$ws.Range('E45').Value = '2026-02-20 07:26:16'
$ws.Range('E46').Value = '2026-02-20 07:26:18'
